# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.173.34"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.520.75"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'536.24"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'139.88"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.525.90"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'0.0996"
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "2.965.95"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "59.126.22"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "2.507.08"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'10.92"
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'4.22"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'321.81"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "'62.56"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").Value = "'0.423"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "'7.76"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D29").Value = "'6.73"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.80"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0766"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "'160.84"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("D35").Value = "'1.45"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("D36").Value = "'18.50"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").Value = "'36.95"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").Value = "'3.64"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "'5.26"
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("D43").Value = "'283.05"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'10.86"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'122.72"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'18.52"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'0.0510"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.0223"
$ws.Range("E51").Value = "  -1.61%  "
